# figure_3.pptx update
#  1. Remove the now-unused "Within/Between-Subjects Treatments" callout
#     shapes (two background rectangles, two labels, one text box and one
#     connector arrow) from the single slide.
#  2. Refresh the auto date field ("datetimeFigureOut") shown on every
#     slide layout and on the slide master from 11/29/22 to 12/6/22.

$p = $ppt.ActivePresentation

# --- 1. Remove the obsolete shapes from slide 1 -----------------------
$s = $p.Slides.Item(1)

$shapesToRemove = @(
    "Rechteck 42",
    "Rechteck 31",
    "Textfeld 32",
    "Textfeld 34",
    "Textfeld 37",
    "Gerade Verbindung mit Pfeil 39"
)

foreach ($shapeName in $shapesToRemove) {
    $s.Shapes.Item($shapeName).Delete()
}

# --- 2. Update the cached date field text ------------------------------
$newDate = "12/6/22"

# Slide master date placeholder
$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Name -like "Datumsplatzhalter*") {
        $shp.TextFrame.TextRange.Text = $newDate
    }
}

# Every slide layout's date placeholder
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Name -like "Datumsplatzhalter*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}
